$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.741.71"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").Value = "2.635.05"
$ws.Range("E3").Value = "  +4.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.79"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.76"
$ws.Range("E6").Value = "  +2.45%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("E8").Value = "  +1.40%  "
$ws.Range("D9").Value = "2.632.85"
$ws.Range("E9").Value = "  +4.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("E10").Value = "  +13.96%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +2.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.06"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").Value = "3.126.09"
$ws.Range("E14").Value = "  +7.47%  "
$ws.Range("E15").Value = "  +7.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.71"
$ws.Range("E16").Value = "  +2.38%  "
$ws.Range("D17").Value = "71.552.79"
$ws.Range("E17").Value = "  +4.33%  "
$ws.Range("D18").Value = "2.632.63"
$ws.Range("E18").Value = "  +5.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "382.31"
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.49"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.88"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("E23").Value = "  +17.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.72"
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.47"
$ws.Range("E25").Value = "  +6.37%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("E27").Value = "  +10.90%  "
$ws.Range("D28").Value = "2.761.47"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "553.36"
$ws.Range("E30").Value = "  +6.73%  "
$ws.Range("D31").Value = "0.0₃0962"
$ws.Range("E31").Value = "  +8.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.11"
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("E33").Value = "  +8.40%  "
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.25"
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.26"
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("E40").Value = "  +7.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  +6.67%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.62"
$ws.Range("E44").Value = "  +10.04%  "
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.53"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.96"
$ws.Range("E47").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.65"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("E50").Value = "  +6.73%  "
$ws.Range("D51").Value = "0.0₆0265"
$ws.Range("E51").Value = "  +5.05%  "
